# Update "想去人数" (number of people interested) figures on the
# "展览" and "全部类型" sheets, as published by the gh-pages data refresh.
$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 8962
    $ws.Range("F4").Value = 447
}
